$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for row 2 (previously row 6's data)
$ws.Range("D2").Value = 45086
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 18000
$ws.Range("Q2").Value = "`$/caja 18 kilos granel"
$ws.Range("R2").Value = "Región del Maule"
$ws.Range("S2").Value = 1000
$ws.Range("T2").Value = 18

# Target values for row 3 (previously row 2's data)
$ws.Range("D3").Value = 44698
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("Q3").Value = "`$/caja 18 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 917
$ws.Range("T3").Value = 18

# Target values for row 4 (previously row 5's data)
$ws.Range("D4").Value = 44334
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12500
$ws.Range("Q4").Value = "`$/caja 12 kilos empedrada"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1042
$ws.Range("T4").Value = 12

# Target values for row 5 (previously row 9's data)
$ws.Range("D5").Value = 44344
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13500
$ws.Range("Q5").Value = "`$/caja 18 kilos granel"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 750
$ws.Range("T5").Value = 18

# Target values for row 6 (previously row 3's data)
$ws.Range("D6").Value = 44316
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 17500
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 17750
$ws.Range("Q6").Value = "`$/caja 16 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1109
$ws.Range("T6").Value = 16

# Target values for row 7 (previously row 4's data)
$ws.Range("D7").Value = 44316
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 16000
$ws.Range("Q7").Value = "`$/caja 16 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 16

# Row 8 is unchanged (maps to itself)

# Target values for row 9 (previously row 7's data)
$ws.Range("D9").Value = 45085
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("Q9").Value = "`$/caja 18 kilos granel"
$ws.Range("R9").Value = "Región del Maule"
$ws.Range("S9").Value = 1028
$ws.Range("T9").Value = 18
